# Updated cryptos list on Fri Oct 18 13:39:56 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while preserving the original cell representation.
# Numeric-looking strings (e.g. "595.86") are coerced to real numbers by a
# plain Range.Value assignment (normal Excel typing behaviour), but the source
# data stores them as literal text -- so for those we briefly force the cell to
# Text format, assign, then clear the format again so no residual style sticks.
function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
Set-Cell "D2" "67.834.64"
Set-Cell "E2" "  +1.56%  "
# Row 3 - Ethereum
Set-Cell "D3" "2.615.86"
Set-Cell "E3" "  +0.63%  "
# Row 4 - TetherUSD
Set-Cell "E4" "  -0.07%  "
# Row 5 - BNB
Set-TextCell "D5" "595.86"
Set-Cell "E5" "  +0.98%  "
# Row 6 - Solana
Set-TextCell "D6" "153.06"
Set-Cell "E6" "  +0.78%  "
# Row 7 - USDC
Set-Cell "E7" "  -0.06%  "
# Row 8 - XRP
Set-TextCell "D8" "0.543"
Set-Cell "E8" "  -1.56%  "
# Row 9 - LidoStakedEther
Set-Cell "D9" "2.615.85"
Set-Cell "E9" "  +0.47%  "
# Row 10 - Dogecoin
Set-Cell "E10" "  +9.96%  "
# Row 11 - TRON
Set-Cell "E11" "  -0.45%  "
# Row 12 - Toncoin
Set-TextCell "D12" "5.22"
Set-Cell "E12" "  +0.85%  "
# Row 13 - Cardano
Set-TextCell "D13" "0.347"
Set-Cell "E13" "  +0.45%  "
# Row 14 - Avalanche
Set-TextCell "D14" "27.55"
Set-Cell "E14" "  +0.40%  "
# Row 15 - ShibaInu
Set-Cell "E15" "  +4.96%  "
# Row 16 - WrappedliquidstakedEther2.0
Set-Cell "D16" "3.094.78"
Set-Cell "E16" "  +0.57%  "
# Row 17 - WrappedBTC
Set-Cell "D17" "67.742.50"
Set-Cell "E17" "  +1.56%  "
# Row 18 - WrappedEther
Set-Cell "D18" "2.617.42"
Set-Cell "E18" "  +0.46%  "
# Row 19 - Chainlink
Set-TextCell "D19" "11.33"
Set-Cell "E19" "  +3.06%  "
# Row 20 - BitcoinCash
Set-TextCell "D20" "366.81"
Set-Cell "E20" "  +0.68%  "
# Row 21 - Uniswap
Set-TextCell "D21" "7.38"
Set-Cell "E21" "  +0.59%  "
# Row 22 - Polkadot
Set-TextCell "D22" "4.21"
Set-Cell "E22" "  -1.73%  "
# Row 23 - NEARProtocol
Set-TextCell "D23" "4.78"
Set-Cell "E23" "  -0.89%  "
# Row 24 - SuiNetwork
Set-TextCell "D24" "2.07"
Set-Cell "E24" "  +1.27%  "
# Row 25 - Litecoin
Set-TextCell "D25" "72.72"
Set-Cell "E25" "  +8.04%  "
# Row 26 - Dai
Set-Cell "E26" "  -0.03%  "
# Row 27 - Aptos
Set-TextCell "D27" "9.91"
Set-Cell "E27" "  -1.31%  "
# Row 29 - PEPE
Set-Cell "E29" "  +3.30%  "
# Row 30 - Binance-PegBSC-USD
Set-Cell "E30" "  +0.15%  "
# Row 31 - Bittensor
Set-TextCell "D31" "567.32"
Set-Cell "E31" "  -2.74%  "
# Row 32 <-> Row 33 swap: row 32 now holds Fetch.AI data
Set-Cell "B32" "Fetch.AI"
Set-Cell "C32" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D32" "1.39"
Set-Cell "E32" "  +1.26%  "
# Row 33 now holds InternetComputer(DFINITY) data
Set-Cell "B33" "InternetComputer(DFINITY)"
Set-Cell "C33" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D33" "7.85"
Set-Cell "E33" "  +2.31%  "
# Row 34 - PancakeSwap
Set-Cell "E34" "  +1.53%  "
# Row 35 - FirstDigitalUSD
Set-TextCell "D35" "0.998"
Set-Cell "E35" "  -0.10%  "
# Row 36 - Kaspa
Set-TextCell "D36" "0.126"
Set-Cell "E36" "  +3.21%  "
# Row 37 - ImmutableX
Set-Cell "E37" "  +1.15%  "
# Row 38 - Monero
Set-TextCell "D38" "161.69"
Set-Cell "E38" "  +4.45%  "
# Row 39 - EthereumClassic
Set-TextCell "D39" "19.06"
Set-Cell "E39" "  +1.06%  "
# Row 40 - Stacks
Set-TextCell "D40" "1.88"
Set-Cell "E40" "  +4.95%  "
# Row 41 - PolygonEcosystemToken
Set-Cell "E41" "  +0.44%  "
# Row 42 - RenderToken
Set-Cell "E42" "  +1.56%  "
# Row 43 - dogwifhat
Set-Cell "E43" "  +2.74%  "
# Row 44 - BabyDogeCoin
Set-Cell "E44" "  +12.59%  "
# Row 45 - WhiteBITCoin
Set-Cell "E45" "  +3.77%  "
# Row 47 - OKB
Set-TextCell "D47" "40.12"
Set-Cell "E47" "  -1.58%  "
# Row 48 - Aave
Set-TextCell "D48" "154.77"
Set-Cell "E48" "  +0.55%  "
# Row 49 - Filecoin
Set-Cell "E49" "  -1.21%  "
# Row 50 - InjectiveProtocol
Set-Cell "E50" "  +1.51%  "
# Row 51 - Optimism
Set-TextCell "D51" "1.68"
Set-Cell "E51" "  -0.15%  "

